$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 3 (product 5414150631147) entirely - this shifts rows 4-8 up to 3-7
$ws.Rows.Item(3).Delete()

# Update the expiry dates (column C) for the remaining rows to their new values
$ws.Range("C2").Value = 44592
$ws.Range("C4").Value = 44592
$ws.Range("C5").Value = 44593
$ws.Range("C7").Value = 44592

# Update the active selection to match the new last cell
$ws.Range("C7").Select()
